$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B30: "Goals" -> " Goals + Use cases"
$ws.Range("B30").Value = " Goals + Use cases"

# Update C30: 1 -> 2
$ws.Range("C30").Value = 2

# Update the selection to match F29
$ws.Range("F29").Select()
